# Correction in SA algorithm and 746 logs
# Update the Fitness column (C) values for run_6.xlsx log sheet.
# Rows 2-33   (Generation 0-31)   -> 7812
# Rows 34-49  (Generation 32-47)  -> 7318
# Rows 50-252 (Generation 48-250) -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C33").Value = 7812
$ws.Range("C34:C49").Value = 7318
$ws.Range("C50:C252").Value = 7293
